$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the two timestamp cells on the existing last row (row 8) ---
# Tiny precision corrections to the already-recorded stop/start timestamps.
$ws.Range("C8").Value = 45406.95036223379
$ws.Range("D8").Value = 45406.95041997685

# --- Append 5 new measurement rows (9-13), copying the number formatting
#     (date/time style for Time start/Time End, time style for Productive
#     time) down from row 8 before filling in the new values. ---
$newRows = 9..13

foreach ($r in $newRows) {
    $destRange = "A{0}:K{0}" -f $r
    $ws.Range("A8:K8").Copy()
    $ws.Range($destRange).PasteSpecial(-4122)
}

# Row 9
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Snakes"
$ws.Range("C9").Value = 45406.96536712963
$ws.Range("D9").Value = 45406.96547386574
$ws.Range("F9").Value = 0.0001041666666666667
$ws.Range("H9").Value = "Team1"
$ws.Range("I9").Value = "Process1"
$ws.Range("J9").Value = "Person1"
$ws.Range("E9").Clear()
$ws.Range("G9").Clear()
$ws.Range("K9").Clear()

# Row 10
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Snakes"
$ws.Range("C10").Value = 45406.96608680556
$ws.Range("D10").Value = 45406.96632180556
$ws.Range("F10").Value = 0.0002314814814814815
$ws.Range("H10").Value = "Team2"
$ws.Range("I10").Value = "Process15"
$ws.Range("J10").Value = "Person9"
$ws.Range("E10").Clear()
$ws.Range("G10").Clear()
$ws.Range("K10").Clear()

# Row 11
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "Snakes"
$ws.Range("C11").Value = 45406.99643114583
$ws.Range("D11").Value = 45406.99703233796
$ws.Range("F11").Value = 0.0005902777777777778
$ws.Range("H11").Value = "Team1"
$ws.Range("I11").Value = "Process1"
$ws.Range("J11").Value = "Person1"
$ws.Range("E11").Clear()
$ws.Range("G11").Clear()
$ws.Range("K11").Clear()

# Row 12
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "Snakes"
$ws.Range("C12").Value = 45406.999475
$ws.Range("D12").Value = 45406.99950525463
$ws.Range("F12").Value = 0.00002314814814814815
$ws.Range("H12").Value = "Team1"
$ws.Range("I12").Value = "Process1"
$ws.Range("J12").Value = "Person1"
$ws.Range("E12").Clear()
$ws.Range("G12").Clear()
$ws.Range("K12").Clear()

# Row 13
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "Snakes"
$ws.Range("C13").Value = 45407.35814317341
$ws.Range("D13").Value = 45407.35835694815
$ws.Range("F13").Value = 0.0002083333333333333
$ws.Range("H13").Value = "Team1"
$ws.Range("I13").Value = "Process1"
$ws.Range("J13").Value = "Person1"
$ws.Range("E13").Clear()
$ws.Range("G13").Clear()
$ws.Range("K13").Clear()
